$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-17 (old Resolving-Mac target rows removed)
$ws.Rows("14:17").Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 7).Value = 189.0573523333333
$ws.Cells.Item(2, 8).Value = 567.172057
$ws.Cells.Item(2, 9).Value = 0.1182556374491171
$ws.Cells.Item(2, 10).Value = 0.1182556374491171
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.4568563333333334
$ws.Cells.Item(2, 14).Value = 1.370569
$ws.Cells.Item(2, 15).Value = 0.2459930547478846
$ws.Cells.Item(2, 16).Value = 0.2459930547478847
$ws.Cells.Item(2, 17).Value = 86.37204877671479
$ws.Cells.Item(2, 18).Value = 777.3484389904331
$ws.Cells.Item(2, 19).Value = 0.02909006549726666
$ws.Cells.Item(2, 20).Value = 0.02909006549726666

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 7).Value = 189.0573523333333
$ws.Cells.Item(3, 8).Value = 567.172057
$ws.Cells.Item(3, 9).Value = 0.1182556374491171
$ws.Cells.Item(3, 10).Value = 0.1182556374491171
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.194006666666667
$ws.Cells.Item(3, 14).Value = 3.58202
$ws.Cells.Item(3, 15).Value = 0.6429096542881224
$ws.Cells.Item(3, 16).Value = 0.6429096542881225
$ws.Cells.Item(3, 17).Value = 225.7357390683489
$ws.Cells.Item(3, 18).Value = 2031.62165161514
$ws.Cells.Item(3, 19).Value = 0.07602769099003341
$ws.Cells.Item(3, 20).Value = 0.07602769099003343

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 7).Value = 189.0573523333333
$ws.Cells.Item(4, 8).Value = 567.172057
$ws.Cells.Item(4, 9).Value = 0.1182556374491171
$ws.Cells.Item(4, 10).Value = 0.1182556374491171
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.206329
$ws.Cells.Item(4, 14).Value = 0.6189870000000001
$ws.Cells.Item(4, 15).Value = 0.111097290963993
$ws.Cells.Item(4, 16).Value = 0.111097290963993
$ws.Cells.Item(4, 17).Value = 39.00801444958434
$ws.Cells.Item(4, 18).Value = 351.072130046259
$ws.Cells.Item(4, 19).Value = 0.01313788096181702
$ws.Cells.Item(4, 20).Value = 0.01313788096181702

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 7).Value = 930.1503093333332
$ws.Cells.Item(5, 8).Value = 2790.450928
$ws.Cells.Item(5, 9).Value = 0.5818103152093762
$ws.Cells.Item(5, 10).Value = 0.5818103152093762
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.4568563333333334
$ws.Cells.Item(5, 14).Value = 1.370569
$ws.Cells.Item(5, 15).Value = 0.2459930547478846
$ws.Cells.Item(5, 16).Value = 0.2459930547478847
$ws.Cells.Item(5, 17).Value = 424.9450597708924
$ws.Cells.Item(5, 18).Value = 3824.505537938032
$ws.Cells.Item(5, 19).Value = 0.1431212967221841
$ws.Cells.Item(5, 20).Value = 0.1431212967221841

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 7).Value = 930.1503093333332
$ws.Cells.Item(6, 8).Value = 2790.450928
$ws.Cells.Item(6, 9).Value = 0.5818103152093762
$ws.Cells.Item(6, 10).Value = 0.5818103152093762
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.194006666666667
$ws.Cells.Item(6, 14).Value = 3.58202
$ws.Cells.Item(6, 15).Value = 0.6429096542881224
$ws.Cells.Item(6, 16).Value = 0.6429096542881225
$ws.Cells.Item(6, 17).Value = 1110.605670346062
$ws.Cells.Item(6, 18).Value = 9995.45103311456
$ws.Cells.Item(6, 19).Value = 0.3740514686125236
$ws.Cells.Item(6, 20).Value = 0.3740514686125236

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 7).Value = 930.1503093333332
$ws.Cells.Item(7, 8).Value = 2790.450928
$ws.Cells.Item(7, 9).Value = 0.5818103152093762
$ws.Cells.Item(7, 10).Value = 0.5818103152093762
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.206329
$ws.Cells.Item(7, 14).Value = 0.6189870000000001
$ws.Cells.Item(7, 15).Value = 0.111097290963993
$ws.Cells.Item(7, 16).Value = 0.111097290963993
$ws.Cells.Item(7, 17).Value = 191.9169831744373
$ws.Cells.Item(7, 18).Value = 1727.252848569936
$ws.Cells.Item(7, 19).Value = 0.06463754987466852
$ws.Cells.Item(7, 20).Value = 0.06463754987466853

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 7).Value = 420.6651306666666
$ws.Cells.Item(8, 8).Value = 1261.995392
$ws.Cells.Item(8, 9).Value = 0.2631266256807295
$ws.Cells.Item(8, 10).Value = 0.2631266256807295
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.4568563333333334
$ws.Cells.Item(8, 14).Value = 1.370569
$ws.Cells.Item(8, 15).Value = 0.2459930547478846
$ws.Cells.Item(8, 16).Value = 0.2459930547478847
$ws.Cells.Item(8, 17).Value = 192.1835291575609
$ws.Cells.Item(8, 18).Value = 1729.651762418048
$ws.Cells.Item(8, 19).Value = 0.06472732243670583
$ws.Cells.Item(8, 20).Value = 0.06472732243670586

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 7).Value = 420.6651306666666
$ws.Cells.Item(9, 8).Value = 1261.995392
$ws.Cells.Item(9, 9).Value = 0.2631266256807295
$ws.Cells.Item(9, 10).Value = 0.2631266256807295
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.194006666666667
$ws.Cells.Item(9, 14).Value = 3.58202
$ws.Cells.Item(9, 15).Value = 0.6429096542881224
$ws.Cells.Item(9, 16).Value = 0.6429096542881225
$ws.Cells.Item(9, 17).Value = 502.2769704502043
$ws.Cells.Item(9, 18).Value = 4520.492734051839
$ws.Cells.Item(9, 19).Value = 0.169166647950398
$ws.Cells.Item(9, 20).Value = 0.169166647950398

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 7).Value = 420.6651306666666
$ws.Cells.Item(10, 8).Value = 1261.995392
$ws.Cells.Item(10, 9).Value = 0.2631266256807295
$ws.Cells.Item(10, 10).Value = 0.2631266256807295
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.206329
$ws.Cells.Item(10, 14).Value = 0.6189870000000001
$ws.Cells.Item(10, 15).Value = 0.111097290963993
$ws.Cells.Item(10, 16).Value = 0.111097290963993
$ws.Cells.Item(10, 17).Value = 86.79541574532266
$ws.Cells.Item(10, 18).Value = 781.1587417079039
$ws.Cells.Item(10, 19).Value = 0.02923265529362567
$ws.Cells.Item(10, 20).Value = 0.02923265529362567

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 7).Value = 58.84466766666667
$ws.Cells.Item(11, 8).Value = 176.534003
$ws.Cells.Item(11, 9).Value = 0.03680742166077718
$ws.Cells.Item(11, 10).Value = 0.03680742166077718
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.4568563333333334
$ws.Cells.Item(11, 14).Value = 1.370569
$ws.Cells.Item(11, 15).Value = 0.2459930547478846
$ws.Cells.Item(11, 16).Value = 0.2459930547478847
$ws.Cells.Item(11, 17).Value = 26.88355910641189
$ws.Cells.Item(11, 18).Value = 241.952031957707
$ws.Cells.Item(11, 19).Value = 0.009054370091728036
$ws.Cells.Item(11, 20).Value = 0.009054370091728037

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 7).Value = 58.84466766666667
$ws.Cells.Item(12, 8).Value = 176.534003
$ws.Cells.Item(12, 9).Value = 0.03680742166077718
$ws.Cells.Item(12, 10).Value = 0.03680742166077718
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.194006666666667
$ws.Cells.Item(12, 14).Value = 3.58202
$ws.Cells.Item(12, 15).Value = 0.6429096542881224
$ws.Cells.Item(12, 16).Value = 0.6429096542881225
$ws.Cells.Item(12, 17).Value = 70.26092549178445
$ws.Cells.Item(12, 18).Value = 632.34832942606
$ws.Cells.Item(12, 19).Value = 0.02366384673516741
$ws.Cells.Item(12, 20).Value = 0.02366384673516741

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 7).Value = 58.84466766666667
$ws.Cells.Item(13, 8).Value = 176.534003
$ws.Cells.Item(13, 9).Value = 0.03680742166077718
$ws.Cells.Item(13, 10).Value = 0.03680742166077718
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.206329
$ws.Cells.Item(13, 14).Value = 0.6189870000000001
$ws.Cells.Item(13, 15).Value = 0.111097290963993
$ws.Cells.Item(13, 16).Value = 0.111097290963993
$ws.Cells.Item(13, 17).Value = 12.14136143499567
$ws.Cells.Item(13, 18).Value = 109.272252914961
$ws.Cells.Item(13, 19).Value = 0.00408920483388174
$ws.Cells.Item(13, 20).Value = 0.00408920483388174
